$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-09-22 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-23 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("515÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "592÷9=", 2) | Out-Null
$d.Content.Find.Execute("424÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "150÷4=", 2) | Out-Null
$d.Content.Find.Execute("837÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "711÷7=", 2) | Out-Null
$d.Content.Find.Execute("823÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "824÷7=", 2) | Out-Null
$d.Content.Find.Execute("530÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "947÷4=", 2) | Out-Null
$d.Content.Find.Execute("153÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "130÷5=", 2) | Out-Null
$d.Content.Find.Execute("729÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "763÷2=", 2) | Out-Null
$d.Content.Find.Execute("856÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "123÷7=", 2) | Out-Null
$d.Content.Find.Execute("319÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "621÷5=", 2) | Out-Null
$d.Content.Find.Execute("806÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "716÷5=", 2) | Out-Null
$d.Content.Find.Execute("962÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "357÷3=", 2) | Out-Null
$d.Content.Find.Execute("393÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "704÷4=", 2) | Out-Null
$d.Content.Find.Execute("761÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "340÷2=", 2) | Out-Null
$d.Content.Find.Execute("307÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "530÷9=", 2) | Out-Null
$d.Content.Find.Execute("511÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "356÷5=", 2) | Out-Null
$d.Content.Find.Execute("497÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "159÷3=", 2) | Out-Null
$d.Content.Find.Execute("845÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "297÷7=", 2) | Out-Null
$d.Content.Find.Execute("661÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "837÷2=", 2) | Out-Null
$d.Content.Find.Execute("232÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "846÷3=", 2) | Out-Null
$d.Content.Find.Execute("280÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "324÷9=", 2) | Out-Null
$d.Content.Find.Execute("715÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "942÷8=", 2) | Out-Null
$d.Content.Find.Execute("132÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "923÷9=", 2) | Out-Null
$d.Content.Find.Execute("233÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "375÷7=", 2) | Out-Null
$d.Content.Find.Execute("995÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "775÷4=", 2) | Out-Null
$d.Content.Find.Execute("408÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "227÷5=", 2) | Out-Null
